# issue #5: stock data from json to db
# Adds "category", "source_file" and "index" columns to the 股票 (stock) sheet,
# and fixes a stray bullet character in a company name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Fix stray bullet character in the shared string text -----------------
# "仁寶電腦工業股份有限公司•" -> "仁寶電腦工業股份有限公司"
if ($ws.Cells.Item(3, 2).Value2 -eq "仁寶電腦工業股份有限公司•") {
    $ws.Cells.Item(3, 2).Value = "仁寶電腦工業股份有限公司"
}

# --- Insert a new "category" column right after "property_category" ------
# Old layout: B name | C owner | D quantity | E face_value | F currency |
#             G total | H property_category | I date | J legislator_name |
#             K legislator_id
# New layout adds "category" (col I) and pushes date/legislator_name/
# legislator_id one column to the right, then appends source_file and index.
$ws.Columns.Item(9).Insert()

$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

$lastRow = 11
for ($r = 2; $r -le $lastRow; $r++) {
    $idx = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmpba991"
    $ws.Cells.Item($r, 14).Value = $idx
}
